$wb = $excel.ActiveWorkbook

# Sheets "展览" (sheetId 1) and "全部类型" (sheetId 4) contain identical data
# and both need the "想去人数" (F column) counts bumped for the same four rows.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F5").Value = 1083
    $ws.Range("F6").Value = 784
    $ws.Range("F7").Value = 37
    $ws.Range("F8").Value = 5813
}
